$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.1993817619783617
$ws.Cells.Item(2, 3).Value = 0.5239567233384853
$ws.Cells.Item(2, 10).Value = 0.02009273570324575
$ws.Cells.Item(2, 15).Value = 0.001545595054095827
$ws.Cells.Item(2, 16).Value = 0.160741885625966
$ws.Cells.Item(2, 19).Value = 0.09428129829984544

# Row 3
$ws.Cells.Item(3, 2).Value = 0.01534526854219949
$ws.Cells.Item(3, 3).Value = 0.02813299232736573
$ws.Cells.Item(3, 10).Value = 0.04603580562659847
$ws.Cells.Item(3, 16).Value = 0.7442455242966752
$ws.Cells.Item(3, 19).Value = 0.1662404092071611

# Row 4
$ws.Cells.Item(4, 10).Value = 0.1089108910891089
$ws.Cells.Item(4, 16).Value = 0.594059405940594
$ws.Cells.Item(4, 19).Value = 0.297029702970297

# Row 6
$ws.Cells.Item(6, 2).Value = 0.0772626931567329
$ws.Cells.Item(6, 4).Value = 0.01324503311258278
$ws.Cells.Item(6, 5).Value = 0.002207505518763797
$ws.Cells.Item(6, 6).Value = 0.05960264900662252
$ws.Cells.Item(6, 10).Value = 0.2759381898454746
$ws.Cells.Item(6, 15).Value = 0.02428256070640177
$ws.Cells.Item(6, 17).Value = 0.1920529801324503
$ws.Cells.Item(6, 18).Value = 0.04194260485651214
$ws.Cells.Item(6, 19).Value = 0.3134657836644592

# Row 7
$ws.Cells.Item(7, 2).Value = 0.1355498721227621
$ws.Cells.Item(7, 4).Value = 0.007672634271099744
$ws.Cells.Item(7, 6).Value = 0.06393861892583121
$ws.Cells.Item(7, 10).Value = 0.1304347826086956
$ws.Cells.Item(7, 15).Value = 0.02813299232736573
$ws.Cells.Item(7, 17).Value = 0.1457800511508952
$ws.Cells.Item(7, 18).Value = 0.07928388746803069
$ws.Cells.Item(7, 19).Value = 0.4092071611253197

# Row 8
$ws.Cells.Item(8, 2).Value = 0.1039260969976905
$ws.Cells.Item(8, 4).Value = 0.02309468822170901
$ws.Cells.Item(8, 5).Value = 0.00115473441108545
$ws.Cells.Item(8, 6).Value = 0.06351039260969978
$ws.Cells.Item(8, 10).Value = 0.1131639722863741
$ws.Cells.Item(8, 15).Value = 0.02655889145496536
$ws.Cells.Item(8, 17).Value = 0.2032332563510393
$ws.Cells.Item(8, 18).Value = 0.07159353348729793
$ws.Cells.Item(8, 19).Value = 0.3937644341801386

# Row 9
$ws.Cells.Item(9, 2).Value = 0.08278145695364239
$ws.Cells.Item(9, 4).Value = 0.01986754966887417
$ws.Cells.Item(9, 6).Value = 0.05629139072847682
$ws.Cells.Item(9, 10).Value = 0.1291390728476821
$ws.Cells.Item(9, 15).Value = 0.03642384105960265
$ws.Cells.Item(9, 17).Value = 0.1754966887417219
$ws.Cells.Item(9, 18).Value = 0.08940397350993377
$ws.Cells.Item(9, 19).Value = 0.4105960264900662

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1113258600695787
$ws.Cells.Item(10, 4).Value = 0.02435253189022033
$ws.Cells.Item(10, 5).Value = 0.001932740626207963
$ws.Cells.Item(10, 6).Value = 0.06494008504058756
$ws.Cells.Item(10, 10).Value = 0.1310398144568999
$ws.Cells.Item(10, 15).Value = 0.02821801314263626
$ws.Cells.Item(10, 17).Value = 0.2365674526478547
$ws.Cells.Item(10, 18).Value = 0.06957866254348666
$ws.Cells.Item(10, 19).Value = 0.332044839582528

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1414634146341463
$ws.Cells.Item(11, 10).Value = 0.1040650406504065
$ws.Cells.Item(11, 11).Value = 0.2016260162601626
$ws.Cells.Item(11, 12).Value = 0.5398373983739837
$ws.Cells.Item(11, 19).Value = 0.01300813008130081

# Row 12
$ws.Cells.Item(12, 7).Value = 0.7048710601719198
$ws.Cells.Item(12, 10).Value = 0.2177650429799427
$ws.Cells.Item(12, 11).Value = 0.01146131805157593
$ws.Cells.Item(12, 12).Value = 0.02865329512893983
$ws.Cells.Item(12, 19).Value = 0.03724928366762178

# Row 13
$ws.Cells.Item(13, 7).Value = 0.7674418604651163
$ws.Cells.Item(13, 10).Value = 0.1279069767441861
$ws.Cells.Item(13, 19).Value = 0.1046511627906977

# Row 14
$ws.Cells.Item(14, 7).Value = 0.5
$ws.Cells.Item(14, 10).Value = 0.25
$ws.Cells.Item(14, 19).Value = 0.25

# Row 15
$ws.Cells.Item(15, 6).Value = 0.03382663847780127
$ws.Cells.Item(15, 8).Value = 0.1585623678646934
$ws.Cells.Item(15, 9).Value = 0.04651162790697674
$ws.Cells.Item(15, 10).Value = 0.3361522198731501
$ws.Cells.Item(15, 11).Value = 0.07399577167019028
$ws.Cells.Item(15, 13).Value = 0.0105708245243129
$ws.Cells.Item(15, 14).Value = 0.004228329809725159
$ws.Cells.Item(15, 15).Value = 0.05285412262156448
$ws.Cells.Item(15, 19).Value = 0.2832980972515856

# Row 16
$ws.Cells.Item(16, 6).Value = 0.02073732718894009
$ws.Cells.Item(16, 8).Value = 0.1728110599078341
$ws.Cells.Item(16, 9).Value = 0.05069124423963134
$ws.Cells.Item(16, 10).Value = 0.4009216589861751
$ws.Cells.Item(16, 11).Value = 0.1129032258064516
$ws.Cells.Item(16, 13).Value = 0.01382488479262673
$ws.Cells.Item(16, 15).Value = 0.04838709677419355
$ws.Cells.Item(16, 19).Value = 0.1797235023041475

# Row 17
$ws.Cells.Item(17, 6).Value = 0.02288557213930348
$ws.Cells.Item(17, 8).Value = 0.1781094527363184
$ws.Cells.Item(17, 9).Value = 0.06766169154228856
$ws.Cells.Item(17, 10).Value = 0.4328358208955224
$ws.Cells.Item(17, 11).Value = 0.09751243781094528
$ws.Cells.Item(17, 13).Value = 0.01691542288557214
$ws.Cells.Item(17, 14).Value = 0.0009950248756218905
$ws.Cells.Item(17, 15).Value = 0.06368159203980099
$ws.Cells.Item(17, 19).Value = 0.1194029850746269

# Row 18
$ws.Cells.Item(18, 6).Value = 0.0182370820668693
$ws.Cells.Item(18, 8).Value = 0.1458966565349544
$ws.Cells.Item(18, 9).Value = 0.0851063829787234
$ws.Cells.Item(18, 10).Value = 0.4224924012158054
$ws.Cells.Item(18, 11).Value = 0.1063829787234043
$ws.Cells.Item(18, 13).Value = 0.0060790273556231
$ws.Cells.Item(18, 14).Value = 0.00303951367781155
$ws.Cells.Item(18, 15).Value = 0.0911854103343465
$ws.Cells.Item(18, 19).Value = 0.121580547112462

# Row 19
$ws.Cells.Item(19, 6).Value = 0.02446225221425559
$ws.Cells.Item(19, 8).Value = 0.2113032475748629
$ws.Cells.Item(19, 9).Value = 0.06832560101223113
$ws.Cells.Item(19, 10).Value = 0.3766343315056938
$ws.Cells.Item(19, 11).Value = 0.1197806832560101
$ws.Cells.Item(19, 13).Value = 0.02446225221425559
$ws.Cells.Item(19, 14).Value = 0.0008435259384226065
$ws.Cells.Item(19, 15).Value = 0.07043441585828764
$ws.Cells.Item(19, 19).Value = 0.1037536904259806

